$wb = $excel.ActiveWorkbook

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 882.8
$ws.Range("I28").Value = 658.64703
$ws.Range("K28").Value = 658.64703
$ws.Range("M28").Value = -173.64703

# ALC row 52
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J52").Value = 200
$ws.Range("L52").Value = 600
$ws.Range("N52").Value = -920

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2216.5833
$ws.Range("I116").Value = 1939.9
$ws.Range("J116").Value = 3600
$ws.Range("K116").Value = 1939.9
$ws.Range("L116").Value = 3600
$ws.Range("M116").Value = 1502.1
$ws.Range("N116").Value = -10484

# ALC row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 71429380
$ws.Range("I127").Value = 619.7143
$ws.Range("K127").Value = 1859.1429
$ws.Range("M127").Value = 3100.8571

# ALC row 136
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 55890.453
$ws.Range("J136").Value = 55890.453
$ws.Range("L136").Value = 55890.453
$ws.Range("N136").Value = -66090.45300000001

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1877.921
$ws.Range("I137").Value = 1865.3667
$ws.Range("J137").Value = 1925
$ws.Range("K137").Value = 5596.1001
$ws.Range("L137").Value = 5775
$ws.Range("M137").Value = -3046.1001
$ws.Range("N137").Value = -10875

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3281.4443
$ws.Range("I2").Value = 3281.4443
$ws.Range("K2").Value = 3281.4443
$ws.Range("M2").Value = -3168.4443

# ARM row 29
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 5000
$ws.Range("J29").Value = 5000
$ws.Range("L29").Value = 5000
$ws.Range("N29").Value = -5616

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8049.467
$ws.Range("I61").Value = 4196.8438
$ws.Range("J61").Value = 17532.846
$ws.Range("K61").Value = 4196.8438
$ws.Range("L61").Value = 17532.846
$ws.Range("M61").Value = -3984.8438
$ws.Range("N61").Value = -17956.846

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3713.7446
$ws.Range("I74").Value = 1437.95
$ws.Range("J74").Value = 16718.285
$ws.Range("K74").Value = 1437.95
$ws.Range("L74").Value = 16718.285
$ws.Range("M74").Value = -563.95
$ws.Range("N74").Value = -18466.285

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3713.7446
$ws.Range("I77").Value = 1437.95
$ws.Range("J77").Value = 16718.285
$ws.Range("K77").Value = 7189.75
$ws.Range("L77").Value = 83591.425
$ws.Range("M77").Value = -2821.75
$ws.Range("N77").Value = -92327.425

# ARM row 82
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H82").Value = 37533.375
$ws.Range("J82").Value = 37533.375
$ws.Range("L82").Value = 37533.375
$ws.Range("N82").Value = -38255.375

# ARM row 85
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H85").Value = 37533.375
$ws.Range("J85").Value = 37533.375
$ws.Range("L85").Value = 37533.375
$ws.Range("N85").Value = -40029.375

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 3281.4443
$ws.Range("I116").Value = 3281.4443
$ws.Range("K116").Value = 3281.4443
$ws.Range("M116").Value = -987.4443000000001

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2366.5715
$ws.Range("I122").Value = 2491.2727
$ws.Range("J122").Value = 1909.3334
$ws.Range("K122").Value = 7473.8181
$ws.Range("L122").Value = 5728.0002
$ws.Range("M122").Value = -5023.8181
$ws.Range("N122").Value = -10628.0002

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 8049.467
$ws.Range("I136").Value = 4196.8438
$ws.Range("J136").Value = 17532.846
$ws.Range("K136").Value = 12590.5314
$ws.Range("L136").Value = 52598.538
$ws.Range("M136").Value = -10040.5314
$ws.Range("N136").Value = -57698.538

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3281.4443
$ws.Range("I3").Value = 3281.4443
$ws.Range("K3").Value = 3281.4443
$ws.Range("M3").Value = -3167.4443

# BSM row 12
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 3000
$ws.Range("J12").Value = 3000
$ws.Range("L12").Value = 3000
$ws.Range("N12").Value = -3336

# BSM row 18
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 622.2222
$ws.Range("I105").Value = 576.25
$ws.Range("K105").Value = 576.25
$ws.Range("M105").Value = 1170.75

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 535.84
$ws.Range("I107").Value = 285.63635
$ws.Range("K107").Value = 285.63635
$ws.Range("M107").Value = 1634.36365

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 5521.4375
$ws.Range("I122").Value = 6559.5835
$ws.Range("J122").Value = 2407
$ws.Range("K122").Value = 19678.7505
$ws.Range("L122").Value = 7221
$ws.Range("M122").Value = -17228.7505
$ws.Range("N122").Value = -12121

# CUL row 20
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 1640.1
$ws.Range("I20").Value = 733.5
$ws.Range("K20").Value = 2200.5
$ws.Range("M20").Value = -1973.5

# CUL row 35
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 1500
$ws.Range("J35").Value = 1500
$ws.Range("L35").Value = 4500
$ws.Range("N35").Value = -5076

# CUL row 63
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 3566.5557
$ws.Range("I63").Value = 2406
$ws.Range("J63").Value = 3898.1428
$ws.Range("K63").Value = 7218
$ws.Range("L63").Value = 11694.4284
$ws.Range("M63").Value = -6469
$ws.Range("N63").Value = -13192.4284

# CUL row 66
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 3566.5557
$ws.Range("I66").Value = 2406
$ws.Range("J66").Value = 3898.1428
$ws.Range("K66").Value = 21654
$ws.Range("L66").Value = 35083.2852
$ws.Range("M66").Value = -17910
$ws.Range("N66").Value = -42571.2852

# CUL row 80
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2454.0908
$ws.Range("J80").Value = 2400.7334
$ws.Range("L80").Value = 7202.2002
$ws.Range("N80").Value = -9074.200199999999

# CUL row 83
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 2454.0908
$ws.Range("J83").Value = 2400.7334
$ws.Range("L83").Value = 21606.6006
$ws.Range("N83").Value = -30966.6006

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1490.8
$ws.Range("I97").Value = 1613.5
$ws.Range("K97").Value = 1613.5
$ws.Range("M97").Value = -1117.5

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2660.0625
$ws.Range("I113").Value = 2712.3845
$ws.Range("J113").Value = 2433.3333
$ws.Range("K113").Value = 2712.3845
$ws.Range("L113").Value = 2433.3333
$ws.Range("M113").Value = -542.3845000000001
$ws.Range("N113").Value = -6773.3333

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5416.385
$ws.Range("I122").Value = 6857
$ws.Range("J122").Value = 2175
$ws.Range("K122").Value = 20571
$ws.Range("L122").Value = 6525
$ws.Range("M122").Value = -18121
$ws.Range("N122").Value = -11425

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5965.9165
$ws.Range("I122").Value = 5662.0645
$ws.Range("J122").Value = 6520
$ws.Range("K122").Value = 16986.1935
$ws.Range("L122").Value = 19560
$ws.Range("M122").Value = -14536.1935
$ws.Range("N122").Value = -24460

# WVR row 32
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 10087.667
$ws.Range("I32").Value = 263
$ws.Range("J32").Value = 15000
$ws.Range("K32").Value = 263
$ws.Range("L32").Value = 15000
$ws.Range("M32").Value = 54
$ws.Range("N32").Value = -15634

# WVR row 46
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 36000
$ws.Range("J46").Value = 36000
$ws.Range("L46").Value = 36000
$ws.Range("N46").Value = -36462

# WVR row 123
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 33333.332
$ws.Range("J123").Value = 33333.332
$ws.Range("L123").Value = 33333.332
$ws.Range("N123").Value = -43133.332

# WVR row 134
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 36000
$ws.Range("J134").Value = 36000
$ws.Range("L134").Value = 108000
$ws.Range("N134").Value = -113070
